# Update the benchmark result table (output query) with refreshed timing values.
# These mirror a re-run of the "sift" vs "compiled" $eq benchmark (see commit message:
# "feat: some var and eq optimizations") that produced new timings for columns B-G,
# rows 2-15 of the "chart" worksheet (the table backed by the Power Query output).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 5271.8104999999996
$ws.Cells.Item(2, 3).Value = 3641.6714999999999
$ws.Cells.Item(2, 4).Value = 305.91320000000002
$ws.Cells.Item(2, 5).Value = 2141.2636000000002
$ws.Cells.Item(2, 6).Value = 2106.4760999999999
$ws.Cells.Item(2, 7).Value = 114.37869999999999

$ws.Cells.Item(3, 2).Value = 6184.6657999999998
$ws.Cells.Item(3, 3).Value = 4540.6220000000003
$ws.Cells.Item(3, 4).Value = 309.24779999999998
$ws.Cells.Item(3, 5).Value = 2117.0403999999999
$ws.Cells.Item(3, 6).Value = 2063.2925
$ws.Cells.Item(3, 7).Value = 112.16679999999999

$ws.Cells.Item(4, 2).Value = 5573.8851999999997
$ws.Cells.Item(4, 3).Value = 3662.4731999999999
$ws.Cells.Item(4, 4).Value = 578.34289999999999
$ws.Cells.Item(4, 5).Value = 3036.2192
$ws.Cells.Item(4, 6).Value = 2547.3344000000002
$ws.Cells.Item(4, 7).Value = 600.49509999999998

$ws.Cells.Item(5, 2).Value = 5434.7365
$ws.Cells.Item(5, 3).Value = 3646.6271999999999
$ws.Cells.Item(5, 4).Value = 436.99930000000001
$ws.Cells.Item(5, 5).Value = 4058.3975999999998
$ws.Cells.Item(5, 6).Value = 3975.0457999999999
$ws.Cells.Item(5, 7).Value = 116.0134

$ws.Cells.Item(6, 2).Value = 6342.6779999999999
$ws.Cells.Item(6, 3).Value = 4531.1466
$ws.Cells.Item(6, 4).Value = 420.53820000000002
$ws.Cells.Item(6, 5).Value = 4028.902
$ws.Cells.Item(6, 6).Value = 3921.9492
$ws.Cells.Item(6, 7).Value = 119.1454

$ws.Cells.Item(7, 2).Value = 5558.6812
$ws.Cells.Item(7, 3).Value = 3636.8434999999999
$ws.Cells.Item(7, 4).Value = 586.24379999999996
$ws.Cells.Item(7, 5).Value = 5148.7785999999996
$ws.Cells.Item(7, 6).Value = 4853.2687999999998
$ws.Cells.Item(7, 7).Value = 413.88589999999999

$ws.Cells.Item(8, 2).Value = 5362.6153000000004
$ws.Cells.Item(8, 3).Value = 3558.3658999999998
$ws.Cells.Item(8, 4).Value = 440.084
$ws.Cells.Item(8, 5).Value = 4074.9331000000002
$ws.Cells.Item(8, 6).Value = 3999.4432000000002
$ws.Cells.Item(8, 7).Value = 126.8472

$ws.Cells.Item(9, 2).Value = 5462.4318999999996
$ws.Cells.Item(9, 3).Value = 3638.9023000000002
$ws.Cells.Item(9, 4).Value = 446.16770000000002
$ws.Cells.Item(9, 5).Value = 4092.1021000000001
$ws.Cells.Item(9, 6).Value = 3997.5133000000001
$ws.Cells.Item(9, 7).Value = 126.16549999999999

$ws.Cells.Item(10, 2).Value = 6177.9260000000004
$ws.Cells.Item(10, 3).Value = 3700.3207000000002
$ws.Cells.Item(10, 4).Value = 1090.9472000000001
$ws.Cells.Item(10, 5).Value = 5554.1737999999996
$ws.Cells.Item(10, 6).Value = 4698.9639999999999
$ws.Cells.Item(10, 7).Value = 917.10619999999994

$ws.Cells.Item(11, 2).Value = 6611.5565999999999
$ws.Cells.Item(11, 3).Value = 4581.0162
$ws.Cells.Item(11, 4).Value = 711.08500000000004
$ws.Cells.Item(11, 5).Value = 4088.8957999999998
$ws.Cells.Item(11, 6).Value = 3944.5664999999999
$ws.Cells.Item(11, 7).Value = 195.20679999999999

$ws.Cells.Item(12, 2).Value = 6427.4387999999999
$ws.Cells.Item(12, 3).Value = 4530.9583000000002
$ws.Cells.Item(12, 4).Value = 533.84169999999995
$ws.Cells.Item(12, 5).Value = 4037.2231999999999
$ws.Cells.Item(12, 6).Value = 3993.9502000000002
$ws.Cells.Item(12, 7).Value = 148.38589999999999

$ws.Cells.Item(13, 2).Value = 5574.9628000000002
$ws.Cells.Item(13, 3).Value = 3587.2265000000002
$ws.Cells.Item(13, 4).Value = 622.99890000000005
$ws.Cells.Item(13, 5).Value = 4082.9263999999998
$ws.Cells.Item(13, 6).Value = 3996.8197
$ws.Cells.Item(13, 7).Value = 170.59049999999999

$ws.Cells.Item(14, 2).Value = 6422.9924000000001
$ws.Cells.Item(14, 3).Value = 4557.6728999999996
$ws.Cells.Item(14, 4).Value = 535.77930000000003
$ws.Cells.Item(14, 5).Value = 7952.1080000000002
$ws.Cells.Item(14, 6).Value = 7769.3141999999998
$ws.Cells.Item(14, 7).Value = 281.54700000000003

$ws.Cells.Item(15, 2).Value = 6545.6286
$ws.Cells.Item(15, 3).Value = 4526.7111999999997
$ws.Cells.Item(15, 4).Value = 618.601
$ws.Cells.Item(15, 5).Value = 8198.2713999999996
$ws.Cells.Item(15, 6).Value = 8079.0136000000002
$ws.Cells.Item(15, 7).Value = 145.7611
